# Refresh the crypto price (column D) and 1h volume-change (column E)
# figures with the latest scrape. Values are plain text in the sheet
# (prices use "." as a thousands separator in spots, e.g. "42.104.07",
# and percentages keep padded whitespace), so we assign them as strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.104.07'
$ws.Range("E2").Value = '  -1.83%  '

$ws.Range("D3").Value = '2.244.93'
$ws.Range("E3").Value = '  -1.80%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '247.77'
$ws.Range("E5").Value = '  -1.72%  '

$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  -2.72%  '

$ws.Range("D7").Value = '76.22'
$ws.Range("E7").Value = '  +3.64%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  -3.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.40'
$ws.Range("E10").Value = '  +6.02%  '

$ws.Range("D11").Value = '0.0957'
$ws.Range("E11").Value = '  -2.16%  '

$ws.Range("D12").Value = '7.15'
$ws.Range("E12").Value = '  -3.51%  '

$ws.Range("E13").Value = '  -2.95%  '

$ws.Range("D14").Value = '2.581.90'
$ws.Range("E14").Value = '  -1.75%  '

$ws.Range("D15").Value = '14.83'
$ws.Range("E15").Value = '  -2.90%  '

$ws.Range("D16").Value = '0.863'
$ws.Range("E16").Value = '  -1.05%  '

$ws.Range("D17").Value = '2.256.66'
$ws.Range("E17").Value = '  -1.61%  '

$ws.Range("D18").Value = '42.010.20'
$ws.Range("E18").Value = '  -1.82%  '

$ws.Range("E19").Value = '  -1.81%  '

$ws.Range("D20").Value = '6.14'
$ws.Range("E20").Value = '  -2.70%  '

$ws.Range("D21").Value = '71.89'
$ws.Range("E21").Value = '  -0.82%  '

$ws.Range("E22").Value = '  +3.82%  '

$ws.Range("D23").Value = '230.91'
$ws.Range("E23").Value = '  -2.41%  '

$ws.Range("E24").Value = '  +0.07%  '

$ws.Range("D25").Value = '11.32'
$ws.Range("E25").Value = '  -1.67%  '

$ws.Range("E26").Value = '  -5.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.30'
$ws.Range("E27").Value = '  -4.80%  '

$ws.Range("D28").Value = '7.25'
$ws.Range("E28").Value = '  +12.01%  '

$ws.Range("E29").Value = '  -1.35%  '

$ws.Range("D30").Value = '169.16'
$ws.Range("E30").Value = '  +1.48%  '

$ws.Range("E31").Value = '  -2.11%  '

$ws.Range("D32").Value = '33.94'
$ws.Range("E32").Value = '  +9.21%  '

$ws.Range("E33").Value = '  +1.05%  '

$ws.Range("E34").Value = '  -4.58%  '

$ws.Range("E35").Value = '  +0.22%  '

$ws.Range("E36").Value = '  -0.94%  '

$ws.Range("D37").Value = '4.91'
$ws.Range("E37").Value = '  +3.20%  '

$ws.Range("D38").Value = '14.27'
$ws.Range("E38").Value = '  -1.07%  '

$ws.Range("D39").Value = '0.0301'
$ws.Range("E39").Value = '  -2.35%  '

$ws.Range("E40").Value = '  +0.42%  '

$ws.Range("E41").Value = '  -6.31%  '

$ws.Range("D42").Value = '113.66'
$ws.Range("E42").Value = '  +13.85%  '

$ws.Range("E43").Value = '  -7.00%  '

$ws.Range("D44").Value = '61.15'
$ws.Range("E44").Value = '  -1.02%  '

$ws.Range("D45").Value = '8.71'
$ws.Range("E45").Value = '  -3.94%  '

$ws.Range("E46").Value = '  -2.84%  '

$ws.Range("D47").Value = '0.997'
$ws.Range("E47").Value = '  -0.33%  '

$ws.Range("E48").Value = '  -2.80%  '

$ws.Range("E49").Value = '  -1.32%  '

$ws.Range("D50").Value = '4.25'
$ws.Range("E50").Value = '  -13.08%  '

$ws.Range("D51").Value = '2.28'
$ws.Range("E51").Value = '  -1.00%  '
